$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = -12.2463
$ws.Range("A3").Value = -22.05050000000001
$ws.Range("A14").Value = -21.8015
$ws.Range("A21").Value = -20.03269999999998
$ws.Range("A23").Value = -20.34559999999997
$ws.Range("A25").Value = -21.79050000000001
$ws.Range("C25").Value = -12.86569999999999
$ws.Range("A26").Value = -21.05419999999997
$ws.Range("C27").Value = -12.78829999999999
$ws.Range("A29").Value = -20.74999999999999
$ws.Range("C31").Value = -13.23820000000001
$ws.Range("C39").Value = -12.61900000000001
$ws.Range("C48").Value = -11.16439999999999
$ws.Range("C51").Value = -11.68760000000001
$ws.Range("C52").Value = -11.2066
$ws.Range("A53").Value = -21.93699999999999
$ws.Range("C55").Value = -13.80679999999999
$ws.Range("C56").Value = -12.41819999999999
$ws.Range("A57").Value = -21.9613
$ws.Range("C57").Value = -12.80829999999999
$ws.Range("A59").Value = -22.08929999999999
$ws.Range("A69").Value = -21.59349999999998
$ws.Range("C73").Value = -12.76290000000001
$ws.Range("A79").Value = -20.34200000000001
$ws.Range("A83").Value = -21.9252
$ws.Range("C89").Value = -10.2673
$ws.Range("C90").Value = -12.3132
$ws.Range("A91").Value = -21.43940000000002
$ws.Range("C92").Value = -10.64149999999999
$ws.Range("A93").Value = -20.78429999999999
